$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '60.731.22'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -1.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.384.44'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -1.84%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '568.33'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.88%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '140.77'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -2.54%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '3.385.18'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.84%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.474'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -0.36%  '
$ws.Range('E10').Value = '  -2.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.123'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.66%  '
$ws.Range('E12').Value = '  +2.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '3.961.34'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.89%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '28.39'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.50%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.125'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.13%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '3.387.30'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -2.30%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.0000170'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -1.83%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '60.819.88'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.42%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.34%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.97'
$ws.Range('D20').Style = 'Normal'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '8.96'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.04%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '383.30'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -1.43%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.559'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -1.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '73.67'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.29%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.998'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.20%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000116'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.96%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.524.70'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -1.88%  '
$ws.Range('E28').Value = '  -1.58%  '
$ws.Range('E29').Value = '  -0.22%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.41'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -2.62%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.98'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -2.19%  '
$ws.Range('E33').Value = '  -1.84%  '
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '23.66'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.57%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '6.93'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.89%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '166.89'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.13%  '
$ws.Range('B38').Value = 'RenzoRestakedETH'
$ws.Range('C38').Value = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '3.415.90'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -1.73%  '
$ws.Range('B39').Value = 'NEARProtocol'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.00'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.35%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.48'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.58%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '27.93'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0774'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -0.99%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.780'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.00%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '41.83'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -1.36%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '4.41'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.68%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.67'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -3.32%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.519.58'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.94%  '
$ws.Range('E49').Value = '  -3.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '23.56'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.58%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '6.83'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.43%  '
